$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Link"
$ws.Range("B2").Value = "https://hi-viewer.web.app/mirador/?manifest=https://hi-ut.github.io/dataset/iiif/collection/nishikie.json"
